# Adds the first riddle ("A Cat, a Parrot, and a Bag of Seed") to the top
# of the document, preceded by the author byline and a blank line.
#
# Built as a single raw WordprocessingML fragment and inserted in one shot
# via Range.InsertXML, which lets us reproduce the exact paragraph/run
# formatting (widowControl/autoSpace* pPr flags, Times New Roman rFonts,
# and the gray/black run colors) byte-for-byte instead of reconstructing it
# property-by-property through the object model.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$riddleText = "A man finds himself on a riverbank with a cat, a parrot and a bag of seed.  " +
              "He needs to transport all three to the other side of the river in his boat.  " +
              "However, the boat has room for only the man himself and one other item (either the cat, parrot or seed).  " +
              "In his absence, the cat could eat the parrot, and the parrot would eat the bag of seed.  " +
              "Show how he can get all the passengers to the other side, without leaving the wrong ones alone together."

$grayRunProps  = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="818181"/></w:rPr>'
$blackRunProps = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/></w:rPr>'

$paragraphs = @(
    "<w:p $wNs><w:r><w:t>Steve Swayngim</w:t></w:r></w:p>",

    "<w:p $wNs/>",

    "<w:p $wNs>" +
        "<w:pPr><w:widowControl w:val=`"0`"/><w:autoSpaceDE w:val=`"0`"/><w:autoSpaceDN w:val=`"0`"/><w:adjustRightInd w:val=`"0`"/>$grayRunProps</w:pPr>" +
        "<w:r>$grayRunProps<w:t>A Cat, a Parrot, and a Bag of Seed:</w:t></w:r>" +
    "</w:p>",

    "<w:p $wNs>" +
        "<w:pPr><w:widowControl w:val=`"0`"/><w:autoSpaceDE w:val=`"0`"/><w:autoSpaceDN w:val=`"0`"/><w:adjustRightInd w:val=`"0`"/>$blackRunProps</w:pPr>" +
        "<w:r>$blackRunProps<w:t>$riddleText</w:t></w:r>" +
    "</w:p>"
)

$xml = [string]::Join("", $paragraphs)

# Insert the fragment as the very first thing in the document body, ahead
# of the existing (empty, bookmarked) paragraph.
$insertionPoint = $d.Range(0, 0)
[void]$insertionPoint.InsertXML($xml)
